$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) figures refreshed by the GitHub Actions bot run.
# Values are stored as literal text in the sheet (e.g. "299.80", "-7.00%"), so
# each updated cell is explicitly marked Text first; otherwise Excel would
# auto-convert the entry into a real number/percentage and silently drop
# formatting such as trailing zeros ("299.80" -> 299.8).
$updates = [ordered]@{
    "D2" = "299.80"
    "E2" = "-7.00%"
    "D3" = "35.10"
    "E3" = "-3.06%"
    "D4" = "4.980"
    "E4" = "-2.79%"
    "D5" = "0.07922"
    "E5" = "-1.84%"
    "D6" = "1.908"
    "E6" = "-11.40%"
    "E7" = "-2.62%"
    "D8" = "7.732"
    "E8" = "-4.15%"
    "E9" = "3.54%"
    "D10" = "0.9215"
    "E10" = "-0.86%"
    "D11" = "0.1114"
    "E11" = "10.16%"
    "D12" = "0.1830"
    "E12" = "-2.85%"
    "D13" = "0.09260"
    "E13" = "0.12%"
    "D14" = "0.03529"
    "E14" = "-1.16%"
    "D15" = "0.09884"
    "E15" = "-0.50%"
    "D16" = "0.001385"
    "E16" = "-3.93%"
    "D17" = "0.005707"
    "E17" = "0.76%"
    "D18" = "3.494"
    "E18" = "1.03%"
    "E19" = "2.03%"
    "E20" = "-1.68%"
    "E21" = "-0.30%"
    "D22" = "0.2398"
    "E22" = "8.80%"
    "D23" = "0.04497"
    "E23" = "-2.35%"
    "E24" = "-2.55%"
    "D25" = "0.004584"
    "E25" = "-3.40%"
    "D26" = "0.0001249"
    "E26" = "-3.94%"
    "E27" = "-6.85%"
    "D39" = "0.01888"
    "E39" = "-3.82%"
    "D40" = "0.04677"
    "E40" = "-6.08%"
    "D41" = "0.007605"
    "E41" = "-2.48%"
    "D42" = "0.009562"
    "E42" = "24.27%"
    "D43" = "0.1321"
    "E43" = "-5.74%"
    "D44" = "0.002118"
    "E44" = "2.49%"
    "D45" = "0.01108"
    "E45" = "-6.04%"
    "D46" = "0.00006008"
    "E46" = "-5.90%"
    "E47" = "-0.09%"
    "E49" = "-31.37%"
    "D50" = "0.00002099"
    "E50" = "-0.09%"
    "D51" = "0.0001999"
    "E51" = "-0.09%"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Restore the default (unstyled) look: these cells had no explicit style
# before the edit, so drop the temporary Text number format again once the
# values are safely stored.
$ws.Range("D2:E51").Style = "Normal"
